$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite data rows 2-19 in date-sorted order (per source diff) and append new row 19
# Row 2
$ws.Range("A2").Value = 6
$ws.Range('B2').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C2').Value = 'Metropolitana'
$ws.Range("D2").Value = 44392
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 100112035
$ws.Range('G2').Value = 'Bruselas (repollito)'
$ws.Range('H2').Value = 'Sin especificar'
$ws.Range('I2').Value = 'Primera'
$ws.Range("J2").Value = 220
$ws.Range("K2").Value = 23000
$ws.Range("L2").Value = 23000
$ws.Range("M2").Value = 23000
$ws.Range('N2').Value = '$/malla 15 kilos'
$ws.Range('O2').Value = 'Provincia de Quillota'
$ws.Range("P2").Value = 1533
$ws.Range("Q2").Value = 15
$ws.Range('R2').Value = 'Hortaliza'

# Row 3
$ws.Range("A3").Value = 6
$ws.Range('B3').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C3').Value = 'Metropolitana'
$ws.Range("D3").Value = 44398
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 100112035
$ws.Range('G3').Value = 'Bruselas (repollito)'
$ws.Range('H3').Value = 'Sin especificar'
$ws.Range('I3').Value = 'Primera'
$ws.Range("J3").Value = 130
$ws.Range("K3").Value = 20000
$ws.Range("L3").Value = 20000
$ws.Range("M3").Value = 20000
$ws.Range('N3').Value = '$/malla 15 kilos'
$ws.Range('O3').Value = 'Provincia de Quillota'
$ws.Range("P3").Value = 1333
$ws.Range("Q3").Value = 15
$ws.Range('R3').Value = 'Hortaliza'

# Row 4
$ws.Range("A4").Value = 6
$ws.Range('B4').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C4').Value = 'Metropolitana'
$ws.Range("D4").Value = 44722
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = 100112035
$ws.Range('G4').Value = 'Bruselas (repollito)'
$ws.Range('H4').Value = 'Sin especificar'
$ws.Range('I4').Value = 'Primera'
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 18000
$ws.Range("L4").Value = 20000
$ws.Range("M4").Value = 18933
$ws.Range('N4').Value = '$/malla 15 kilos'
$ws.Range('O4').Value = 'Provincia de Quillota'
$ws.Range("P4").Value = 1262
$ws.Range("Q4").Value = 15
$ws.Range('R4').Value = 'Hortaliza'

# Row 5
$ws.Range("A5").Value = 6
$ws.Range('B5').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C5').Value = 'Metropolitana'
$ws.Range("D5").Value = 44400
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = 100112035
$ws.Range('G5').Value = 'Bruselas (repollito)'
$ws.Range('H5').Value = 'Sin especificar'
$ws.Range('I5').Value = 'Primera'
$ws.Range("J5").Value = 130
$ws.Range("K5").Value = 24000
$ws.Range("L5").Value = 24000
$ws.Range("M5").Value = 24000
$ws.Range('N5').Value = '$/malla 15 kilos'
$ws.Range('O5').Value = 'Provincia de Quillota'
$ws.Range("P5").Value = 1600
$ws.Range("Q5").Value = 15
$ws.Range('R5').Value = 'Hortaliza'

# Row 6
$ws.Range("A6").Value = 6
$ws.Range('B6').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C6').Value = 'Metropolitana'
$ws.Range("D6").Value = 44483
$ws.Range("E6").Value = 13
$ws.Range("F6").Value = 100112035
$ws.Range('G6').Value = 'Bruselas (repollito)'
$ws.Range('H6').Value = 'Sin especificar'
$ws.Range('I6').Value = 'Primera'
$ws.Range("J6").Value = 220
$ws.Range("K6").Value = 18000
$ws.Range("L6").Value = 20000
$ws.Range("M6").Value = 18909
$ws.Range('N6').Value = '$/malla 15 kilos'
$ws.Range('O6').Value = 'Provincia de Quillota'
$ws.Range("P6").Value = 1261
$ws.Range("Q6").Value = 15
$ws.Range('R6').Value = 'Hortaliza'

# Row 7
$ws.Range("A7").Value = 6
$ws.Range('B7').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C7').Value = 'Metropolitana'
$ws.Range("D7").Value = 44446
$ws.Range("E7").Value = 13
$ws.Range("F7").Value = 100112035
$ws.Range('G7').Value = 'Bruselas (repollito)'
$ws.Range('H7').Value = 'Sin especificar'
$ws.Range('I7').Value = 'Primera'
$ws.Range("J7").Value = 150
$ws.Range("K7").Value = 22000
$ws.Range("L7").Value = 24000
$ws.Range("M7").Value = 22667
$ws.Range('N7').Value = '$/malla 15 kilos'
$ws.Range('O7').Value = 'Provincia de Quillota'
$ws.Range("P7").Value = 1511
$ws.Range("Q7").Value = 15
$ws.Range('R7').Value = 'Hortaliza'

# Row 8
$ws.Range("A8").Value = 6
$ws.Range('B8').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C8').Value = 'Metropolitana'
$ws.Range("D8").Value = 44391
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = 100112035
$ws.Range('G8').Value = 'Bruselas (repollito)'
$ws.Range('H8').Value = 'Sin especificar'
$ws.Range('I8').Value = 'Primera'
$ws.Range("J8").Value = 160
$ws.Range("K8").Value = 20000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 20000
$ws.Range('N8').Value = '$/malla 15 kilos'
$ws.Range('O8').Value = 'Provincia de Quillota'
$ws.Range("P8").Value = 1333
$ws.Range("Q8").Value = 15
$ws.Range('R8').Value = 'Hortaliza'

# Row 9
$ws.Range("A9").Value = 6
$ws.Range('B9').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C9').Value = 'Metropolitana'
$ws.Range("D9").Value = 44399
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 100112035
$ws.Range('G9').Value = 'Bruselas (repollito)'
$ws.Range('H9').Value = 'Sin especificar'
$ws.Range('I9').Value = 'Primera'
$ws.Range("J9").Value = 150
$ws.Range("K9").Value = 22000
$ws.Range("L9").Value = 22000
$ws.Range("M9").Value = 22000
$ws.Range('N9').Value = '$/malla 15 kilos'
$ws.Range('O9').Value = 'Provincia de Quillota'
$ws.Range("P9").Value = 1467
$ws.Range("Q9").Value = 15
$ws.Range('R9').Value = 'Hortaliza'

# Row 10
$ws.Range("A10").Value = 6
$ws.Range('B10').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C10').Value = 'Metropolitana'
$ws.Range("D10").Value = 44365
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = 100112035
$ws.Range('G10').Value = 'Bruselas (repollito)'
$ws.Range('H10').Value = 'Sin especificar'
$ws.Range('I10').Value = 'Primera'
$ws.Range("J10").Value = 580
$ws.Range("K10").Value = 20000
$ws.Range("L10").Value = 22000
$ws.Range("M10").Value = 21103
$ws.Range('N10').Value = '$/malla 15 kilos'
$ws.Range('O10').Value = 'Provincia de Quillota'
$ws.Range("P10").Value = 1407
$ws.Range("Q10").Value = 15
$ws.Range('R10').Value = 'Hortaliza'

# Row 11
$ws.Range("A11").Value = 6
$ws.Range('B11').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C11').Value = 'Metropolitana'
$ws.Range("D11").Value = 44714
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = 100112035
$ws.Range('G11').Value = 'Bruselas (repollito)'
$ws.Range('H11').Value = 'Sin especificar'
$ws.Range('I11').Value = 'Primera'
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 16000
$ws.Range("L11").Value = 17000
$ws.Range("M11").Value = 16400
$ws.Range('N11').Value = '$/malla 15 kilos'
$ws.Range('O11').Value = 'Provincia de Quillota'
$ws.Range("P11").Value = 1093
$ws.Range("Q11").Value = 15
$ws.Range('R11').Value = 'Hortaliza'

# Row 12
$ws.Range("A12").Value = 6
$ws.Range('B12').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C12').Value = 'Metropolitana'
$ws.Range("D12").Value = 44727
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = 100112035
$ws.Range('G12').Value = 'Bruselas (repollito)'
$ws.Range('H12').Value = 'Sin especificar'
$ws.Range('I12').Value = 'Primera'
$ws.Range("J12").Value = 220
$ws.Range("K12").Value = 16000
$ws.Range("L12").Value = 18000
$ws.Range("M12").Value = 16909
$ws.Range('N12').Value = '$/malla 15 kilos'
$ws.Range('O12').Value = 'Provincia de Quillota'
$ws.Range("P12").Value = 1127
$ws.Range("Q12").Value = 15
$ws.Range('R12').Value = 'Hortaliza'

# Row 13
$ws.Range("A13").Value = 6
$ws.Range('B13').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C13').Value = 'Metropolitana'
$ws.Range("D13").Value = 44699
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 100112035
$ws.Range('G13').Value = 'Bruselas (repollito)'
$ws.Range('H13').Value = 'Sin especificar'
$ws.Range('I13').Value = 'Primera'
$ws.Range("J13").Value = 150
$ws.Range("K13").Value = 18000
$ws.Range("L13").Value = 20000
$ws.Range("M13").Value = 18667
$ws.Range('N13').Value = '$/malla 15 kilos'
$ws.Range('O13').Value = 'Provincia de Quillota'
$ws.Range("P13").Value = 1244
$ws.Range("Q13").Value = 15
$ws.Range('R13').Value = 'Hortaliza'

# Row 14
$ws.Range("A14").Value = 6
$ws.Range('B14').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C14').Value = 'Metropolitana'
$ws.Range("D14").Value = 44453
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 100112035
$ws.Range('G14').Value = 'Bruselas (repollito)'
$ws.Range('H14').Value = 'Sin especificar'
$ws.Range('I14').Value = 'Primera'
$ws.Range("J14").Value = 280
$ws.Range("K14").Value = 20000
$ws.Range("L14").Value = 22000
$ws.Range("M14").Value = 21286
$ws.Range('N14').Value = '$/malla 15 kilos'
$ws.Range('O14').Value = 'Provincia de Quillota'
$ws.Range("P14").Value = 1419
$ws.Range("Q14").Value = 15
$ws.Range('R14').Value = 'Hortaliza'

# Row 15
$ws.Range("A15").Value = 6
$ws.Range('B15').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C15').Value = 'Metropolitana'
$ws.Range("D15").Value = 44406
$ws.Range("E15").Value = 13
$ws.Range("F15").Value = 100112035
$ws.Range('G15').Value = 'Bruselas (repollito)'
$ws.Range('H15').Value = 'Sin especificar'
$ws.Range('I15').Value = 'Primera'
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 20000
$ws.Range("L15").Value = 22000
$ws.Range("M15").Value = 20850
$ws.Range('N15').Value = '$/malla 15 kilos'
$ws.Range('O15').Value = 'Provincia de Quillota'
$ws.Range("P15").Value = 1390
$ws.Range("Q15").Value = 15
$ws.Range('R15').Value = 'Hortaliza'

# Row 16
$ws.Range("A16").Value = 6
$ws.Range('B16').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C16').Value = 'Metropolitana'
$ws.Range("D16").Value = 44435
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = 100112035
$ws.Range('G16').Value = 'Bruselas (repollito)'
$ws.Range('H16').Value = 'Sin especificar'
$ws.Range('I16').Value = 'Primera'
$ws.Range("J16").Value = 140
$ws.Range("K16").Value = 21000
$ws.Range("L16").Value = 23000
$ws.Range("M16").Value = 21714
$ws.Range('N16').Value = '$/malla 15 kilos'
$ws.Range('O16').Value = 'Provincia de Quillota'
$ws.Range("P16").Value = 1448
$ws.Range("Q16").Value = 15
$ws.Range('R16').Value = 'Hortaliza'

# Row 17
$ws.Range("A17").Value = 6
$ws.Range('B17').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C17').Value = 'Metropolitana'
$ws.Range("D17").Value = 44476
$ws.Range("E17").Value = 13
$ws.Range("F17").Value = 100112035
$ws.Range('G17').Value = 'Bruselas (repollito)'
$ws.Range('H17').Value = 'Sin especificar'
$ws.Range('I17').Value = 'Primera'
$ws.Range("J17").Value = 220
$ws.Range("K17").Value = 20000
$ws.Range("L17").Value = 22000
$ws.Range("M17").Value = 20909
$ws.Range('N17').Value = '$/malla 15 kilos'
$ws.Range('O17').Value = 'Provincia de Quillota'
$ws.Range("P17").Value = 1394
$ws.Range("Q17").Value = 15
$ws.Range('R17').Value = 'Hortaliza'

# Row 18
$ws.Range("A18").Value = 6
$ws.Range('B18').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C18').Value = 'Metropolitana'
$ws.Range("D18").Value = 44449
$ws.Range("E18").Value = 13
$ws.Range("F18").Value = 100112035
$ws.Range('G18').Value = 'Bruselas (repollito)'
$ws.Range('H18').Value = 'Sin especificar'
$ws.Range('I18').Value = 'Primera'
$ws.Range("J18").Value = 220
$ws.Range("K18").Value = 22000
$ws.Range("L18").Value = 24000
$ws.Range("M18").Value = 23091
$ws.Range('N18').Value = '$/malla 15 kilos'
$ws.Range('O18').Value = 'Provincia de Quillota'
$ws.Range("P18").Value = 1539
$ws.Range("Q18").Value = 15
$ws.Range('R18').Value = 'Hortaliza'

# Row 19
$ws.Range("A19").Value = 6
$ws.Range('B19').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C19').Value = 'Metropolitana'
$ws.Range("D19").Value = 44396
$ws.Range("E19").Value = 13
$ws.Range("F19").Value = 100112035
$ws.Range('G19').Value = 'Bruselas (repollito)'
$ws.Range('H19').Value = 'Sin especificar'
$ws.Range('I19').Value = 'Primera'
$ws.Range("J19").Value = 130
$ws.Range("K19").Value = 22000
$ws.Range("L19").Value = 22000
$ws.Range("M19").Value = 22000
$ws.Range('N19').Value = '$/malla 15 kilos'
$ws.Range('O19').Value = 'Provincia de Quillota'
$ws.Range("P19").Value = 1467
$ws.Range("Q19").Value = 15
$ws.Range('R19').Value = 'Hortaliza'

# Match date-column (D) formatting used by the rest of the column for the newly added row
$ws.Range("D19").NumberFormat = $ws.Range("D2").NumberFormat
